$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price-tracker scrape run for 2026-02-07: append one new row to the
# bottom of the table (row 38). Every column in this sheet holds plain
# text (even the numeric-looking price/discount values), so the new
# values are entered as formulas that yield a text result and then
# flattened to static values via copy / paste-special. That keeps the
# stored cells as shared-string text (t="s") without picking up a
# number format / style, exactly like the existing rows.
$row = 38

$ws.Cells.Item($row, 1).Formula = '="2026-02-07"'
$ws.Cells.Item($row, 2).Formula = '="23500000"'
$ws.Cells.Item($row, 3).Formula = '="0"'
$ws.Cells.Item($row, 4).Formula = '="0"'

$rng = $ws.Range("A" + $row + ":D" + $row)
$rng.Copy()
$rng.PasteSpecial(-4163)
